# 12/09/2022 update: new photos - fold the "Seriola rivoliana" rows from
# Feuil2 into the bottom of the Tableau1 data on Feuil1, then drop Feuil2
# (and its now-duplicate "Anisotremus surinamensis" rows) entirely.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws2 = $wb.Worksheets.Item("Feuil2")

# Copy Feuil2's first species (Seriola rivoliana, rows 1-2 = F/M) onto the
# end of Feuil1's table, preserving the row formatting (borders/fills).
$ws2.Range("A1:D2").Copy($ws1.Range("A70")) | Out-Null

# The second species on Feuil2 (Anisotremus surinamensis, rows 3-4) was a
# duplicate entry that is simply being dropped along with the whole sheet.
$ws2.Delete() | Out-Null

# Resync the table/ListObject to the new data extent.
$lo = $ws1.ListObjects.Item("Tableau1")
$lo.Resize($ws1.Range("A1:D71"))

# Restore the sheet's view: scrolled so row 40 is at the top, with A74
# selected just below the new data.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("A74").Select() | Out-Null
